$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report volume/number + week-of dates).
# Each shared string's runs all carry identical formatting in this workbook,
# so replacing the whole cell text reproduces the same visible content.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# ---------------------------------------------------------------------------
# Helper donor cells used to clone number/text formats onto cells that must
# switch between the "-" (text s14) placeholder style and a numeric style.
#   I14 -> numeric style s15 (e.g. used by C/D/F/G integer columns)
#   K14 -> numeric style s16 (e.g. used by E/H/K/L/M/N percent columns)
#   C14 -> text style s14 holding shared string "0"
#   E14 -> text style s14 holding shared string "***.*"
# Copying these into the target cell first clones the correct style + base
# value, then we overwrite with the real value/text we need.
# ---------------------------------------------------------------------------

# Row 15
$ws.Range("F15").Value = 1

# Row 16 (C16 switches from placeholder text "0" to a real numeric count)
$ws.Range("I14").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = 7.142857142857
$ws.Range("L16").Value = 87.5
$ws.Range("M16").Value = -26.229508196721
$ws.Range("N16").Value = -80.349344978165

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 76
$ws.Range("J17").Value = 64
$ws.Range("K17").Value = 18.75
$ws.Range("L17").Value = 40.740740740740
$ws.Range("M17").Value = 5.555555555555
$ws.Range("N17").Value = -64.319248826291

# Row 18
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = -14.634146341463
$ws.Range("L18").Value = 55.555555555555
$ws.Range("M18").Value = -5.405405405405
$ws.Range("N18").Value = -75.609756097561

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 20
$ws.Range("F19").Value = 15
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 111
$ws.Range("J19").Value = 136
$ws.Range("K19").Value = -18.382352941176
$ws.Range("L19").Value = 35.365853658536
$ws.Range("M19").Value = -26.490066225165
$ws.Range("N19").Value = -25.503355704698

# Row 20 (D20/E20 switch from placeholder text to real numeric values)
$ws.Range("I14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 3
$ws.Range("K14").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 45
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = 21.621621621621
$ws.Range("L20").Value = 36.363636363636
$ws.Range("M20").Value = 15.384615384615
$ws.Range("N20").Value = -83.271375464684

# Row 21
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 7.142857142857
$ws.Range("F21").Value = 45
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = -18.181818181818
$ws.Range("I21").Value = 352
$ws.Range("J21").Value = 367
$ws.Range("K21").Value = -4.087193460490
$ws.Range("L21").Value = 44.855967078189
$ws.Range("M21").Value = -12
$ws.Range("N21").Value = -69.681309216192

# Row 23
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 250
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -6.666666666666
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 82
$ws.Range("K23").Value = -2.439024390243
$ws.Range("L23").Value = 42.857142857142
$ws.Range("M23").Value = 42.857142857142

# Row 24
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -46.666666666666
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 23.913043478260
$ws.Range("I24").Value = 406
$ws.Range("J24").Value = 299
$ws.Range("K24").Value = 35.785953177257
$ws.Range("L24").Value = 65.714285714285
$ws.Range("M24").Value = 34.437086092715

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 14
$ws.Range("H25").Value = -22.222222222222
$ws.Range("I25").Value = 109
$ws.Range("J25").Value = 126
$ws.Range("K25").Value = -13.492063492063
$ws.Range("L25").Value = 1.869158878504
$ws.Range("M25").Value = -42.021276595744

# Row 26
$ws.Range("F26").Value = 1

# Row 28 (D28/E28/F28 switch from real numeric values back to placeholder text)
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("F28"))
$ws.Range("H28").Value = -100
$ws.Range("N28").Value = -88.095238095238

# Row 29 (D29/E29/F29 switch from real numeric values back to placeholder text)
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -85.294117647058
